$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Repository"
$ws.Range("B1").Value = "PR Number"
$ws.Range("C1").Value = "Title"
$ws.Range("D1").Value = "Owner"
$ws.Range("E1").Value = "URL"
$ws.Range("F1").Value = "Status"

$ws.Range("A2").Value = "vmn_ecomm_backend"
$ws.Range("B2").Value = 7
$ws.Range("C2").Value = "Model refactoring"
$ws.Range("D2").Value = "jaziel1974"
$ws.Range("E2").Value = "https://github.com/jaziel1974/vmn_ecomm_backend/pull/7"
$ws.Range("F2").Value = "Needs Review"
